$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for data rows 2-130 is a date serial that was
# bumped by one day (46081 -> 46082), i.e. 2026-02-28 -> 2026-03-01.
$ws.Range("C2:C130").Value = 46082
